{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, the\n// \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that separates\n// them from the \"LOM3260: ... (Requisito)\" paragraph above them.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"LOM3260: Computa\u00e7\u00e3o Cient\u00edfica em Python (Requisito)\" paragraph;\n// the three paragraphs that immediately follow it are the ones to remove\n// (a blank paragraph, the \"Ver no Jupiter...\" line, and the \"\u00a9 2020 ...\" line).\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOM3260\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const toRemove = [];\n  for (let i = anchorIndex + 1; i < items.length; i++) {\n    const text = items[i].text.trim();\n    if (text === \"\" || text.indexOf(\"Ver no Jupiter\") !== -1 || text.indexOf(\"Powered by Jekyll\") !== -1) {\n      toRemove.push(items[i]);\n      // Stop once we've captured the blank line + the two footer lines.\n      if (text.indexOf(\"Powered by Jekyll\") !== -1) {\n        break;\n      }\n    } else {\n      break;\n    }\n  }\n  toRemove.forEach((p) => p.delete());\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, the\n# \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that separates\n# them from the \"LOM3260: ... (Requisito)\" paragraph above them.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*LOM3260*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ne -1) {\n    # The paragraph right after the anchor is, in order: a blank paragraph,\n    # the \"Ver no Jupiter...\" line, and the \"\u00a9 2020 ...\" line. Removing the\n    # paragraph that now sits at $anchorIndex + 1 three times in a row deletes\n    # all three, since each deletion shifts the following paragraph into that\n    # slot.\n    for ($n = 0; $n -lt 3; $n++) {\n        $targetIndex = $anchorIndex + 1\n        if ($targetIndex -le $d.Paragraphs.Count) {\n            $d.Paragraphs.Item($targetIndex).Range.Delete()\n        }\n    }\n}\n"}
